$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.456061482429504
$ws.Range("B1").Value = 1.906810879707336
$ws.Range("C1").Value = 3.286095142364502
$ws.Range("D1").Value = 1.389496445655823
$ws.Range("E1").Value = 0.811001718044281
